$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text updated for both language rows ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: fill in Latest Target File / Latest Handback File columns ---
# and record the handback timestamp now that the report has been generated.
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/78606c22659be8a345b553141eb0a5945db356f9/e2e/367c9376-8963-445b-b3c7-d50595daf71e.md", "", "", "367c9376-8963-445b-b3c7-d50595daf71e.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1c6b758d59c28c30801c5d04867d3b7c8aeb6af/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/367c9376-8963-445b-b3c7-d50595daf71e.0acabaea55c351a6c86521a3561361e8b16036ee.zh-cn.xlf", "", "", "367c9376-8963-445b-b3c7-d50595daf71e.0acabaea55c351a6c86521a3561361e8b16036ee.zh-cn.xlf")
$zhcn.Range("H2").Value = "2016-03-14 02:30:49"

$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/78606c22659be8a345b553141eb0a5945db356f9/e2e/e8579b9a-7171-4440-a858-0c2b09a8be45.md", "", "", "e8579b9a-7171-4440-a858-0c2b09a8be45.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1c6b758d59c28c30801c5d04867d3b7c8aeb6af/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e8579b9a-7171-4440-a858-0c2b09a8be45.40aed40c780dac005fd8e5ff0d93493948a1b490.zh-cn.xlf", "", "", "e8579b9a-7171-4440-a858-0c2b09a8be45.40aed40c780dac005fd8e5ff0d93493948a1b490.zh-cn.xlf")
$zhcn.Range("H3").Value = "2016-03-14 02:30:49"

# --- de-de sheet: same treatment ---
$dede = $wb.Worksheets.Item("de-de")

$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/78606c22659be8a345b553141eb0a5945db356f9/e2e/367c9376-8963-445b-b3c7-d50595daf71e.md", "", "", "367c9376-8963-445b-b3c7-d50595daf71e.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6a0f69a52c50eab9dd68058fc8f37e4020f5ac18/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/367c9376-8963-445b-b3c7-d50595daf71e.0acabaea55c351a6c86521a3561361e8b16036ee.de-de.xlf", "", "", "367c9376-8963-445b-b3c7-d50595daf71e.0acabaea55c351a6c86521a3561361e8b16036ee.de-de.xlf")
$dede.Range("H2").Value = "2016-03-14 02:30:54"

$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/78606c22659be8a345b553141eb0a5945db356f9/e2e/e8579b9a-7171-4440-a858-0c2b09a8be45.md", "", "", "e8579b9a-7171-4440-a858-0c2b09a8be45.md")
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6a0f69a52c50eab9dd68058fc8f37e4020f5ac18/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e8579b9a-7171-4440-a858-0c2b09a8be45.40aed40c780dac005fd8e5ff0d93493948a1b490.de-de.xlf", "", "", "e8579b9a-7171-4440-a858-0c2b09a8be45.40aed40c780dac005fd8e5ff0d93493948a1b490.de-de.xlf")
$dede.Range("H3").Value = "2016-03-14 02:30:54"
